# Auto-generated edit script applying the scheduled-runner updates
# to Sheets/Malboro_Profits.xlsx (profit/loss recalculated figures).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 15000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -14532
# Row 23
$ws.Range("H23").Value = 15000
$ws.Range("I23").Value = 15000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 15000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -14766
# Row 29
$ws.Range("H29").Value = 1815
# Row 38
$ws.Range("H38").Value = 1668.125
# Row 40
$ws.Range("H40").Value = 2333.3333
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2325
# Row 62
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376
# Row 65
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880
# Row 70
$ws.Range("H70").Value = 3430.818
$ws.Range("I70").Value = 3707.7778
$ws.Range("J70").Value = 2184.5
$ws.Range("K70").Value = 11123.3334
$ws.Range("L70").Value = 6553.5
$ws.Range("M70").Value = -10853.3334
$ws.Range("N70").Value = -7093.5
# Row 73
$ws.Range("H73").Value = 3430.818
$ws.Range("I73").Value = 3707.7778
$ws.Range("J73").Value = 2184.5
$ws.Range("K73").Value = 11123.3334
$ws.Range("L73").Value = 6553.5
$ws.Range("M73").Value = -10187.3334
$ws.Range("N73").Value = -8425.5
# Row 88
$ws.Range("H88").Value = 4405.357
$ws.Range("I88").Value = 5174.8
$ws.Range("J88").Value = 3977.889
$ws.Range("K88").Value = 5174.8
$ws.Range("L88").Value = 3977.889
$ws.Range("M88").Value = -4768.8
$ws.Range("N88").Value = -4789.889
# Row 91
$ws.Range("H91").Value = 4405.357
$ws.Range("I91").Value = 5174.8
$ws.Range("J91").Value = 3977.889
$ws.Range("K91").Value = 5174.8
$ws.Range("L91").Value = 3977.889
$ws.Range("M91").Value = -3770.8
$ws.Range("N91").Value = -6785.889
# Row 96
$ws.Range("H96").Value = 1882.3914
$ws.Range("I96").Value = 1808.5625
$ws.Range("J96").Value = 2051.1428
$ws.Range("K96").Value = 5425.6875
$ws.Range("L96").Value = 6153.428400000001
$ws.Range("M96").Value = -4052.6875
$ws.Range("N96").Value = -8899.428400000001
# Row 106
$ws.Range("H106").Value = 5509.778
$ws.Range("J106").Value = 3199
$ws.Range("L106").Value = 3199
$ws.Range("N106").Value = -4461
# Row 137
$ws.Range("H137").Value = 19765.084
$ws.Range("J137").Value = 37417.668
$ws.Range("L137").Value = 112253.004
$ws.Range("N137").Value = -117353.004
# Row 138
$ws.Range("H138").Value = 2324.914
$ws.Range("I138").Value = 1755.68
$ws.Range("J138").Value = 2534.1912
$ws.Range("K138").Value = 5267.04
$ws.Range("L138").Value = 7602.573600000001
$ws.Range("M138").Value = -127.04
$ws.Range("N138").Value = -17882.5736
# Row 141
$ws.Range("H141").Value = 3782.7368
$ws.Range("I141").Value = 4254.9287
$ws.Range("J141").Value = 2460.6
$ws.Range("K141").Value = 12764.7861
$ws.Range("L141").Value = 7381.799999999999
$ws.Range("M141").Value = -7584.786100000001
$ws.Range("N141").Value = -17741.8
# Cells removed entirely in the updated workbook
$ws.Range("N21").ClearContents()
$ws.Range("N23").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4287.4375
$ws.Range("I2").Value = 1619.8
$ws.Range("J2").Value = 5500
$ws.Range("K2").Value = 1619.8
$ws.Range("L2").Value = 5500
$ws.Range("M2").Value = -1506.8
$ws.Range("N2").Value = -5726
# Row 32
$ws.Range("H32").Value = 3078.6316
$ws.Range("I32").Value = 1884.0588
$ws.Range("J32").Value = 13232.5
$ws.Range("K32").Value = 1884.0588
$ws.Range("L32").Value = 13232.5
$ws.Range("M32").Value = -1597.0588
$ws.Range("N32").Value = -13806.5
# Row 43
$ws.Range("H43").Value = 33666
$ws.Range("I43").Value = 31998
$ws.Range("K43").Value = 31998
$ws.Range("M43").Value = -31685
# Row 45
$ws.Range("H45").Value = 3091.04
$ws.Range("I45").Value = 3009.5557
$ws.Range("K45").Value = 3009.5557
$ws.Range("M45").Value = -2632.5557
# Row 46
$ws.Range("H46").Value = 5238.1
$ws.Range("J46").Value = 5697.625
$ws.Range("L46").Value = 5697.625
$ws.Range("N46").Value = -6335.625
# Row 116
$ws.Range("H116").Value = 4287.4375
$ws.Range("I116").Value = 1619.8
$ws.Range("J116").Value = 5500
$ws.Range("K116").Value = 1619.8
$ws.Range("L116").Value = 5500
$ws.Range("M116").Value = 674.2
$ws.Range("N116").Value = -10088

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4287.4375
$ws.Range("I3").Value = 1619.8
$ws.Range("J3").Value = 5500
$ws.Range("K3").Value = 1619.8
$ws.Range("L3").Value = 5500
$ws.Range("M3").Value = -1505.8
$ws.Range("N3").Value = -5728

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 80
$ws.Range("H80").Value = 14000
$ws.Range("J80").Value = 14000
$ws.Range("L80").Value = 14000
$ws.Range("N80").Value = -16246
# Row 83
$ws.Range("H83").Value = 14000
$ws.Range("J83").Value = 14000
$ws.Range("L83").Value = 42000
$ws.Range("N83").Value = -53232
# Row 86
$ws.Range("H86").Value = 9716.817999999999
$ws.Range("I86").Value = 11276.125
$ws.Range("K86").Value = 11276.125
$ws.Range("M86").Value = -10153.125
# Row 89
$ws.Range("H89").Value = 9716.817999999999
$ws.Range("I89").Value = 11276.125
$ws.Range("K89").Value = 56380.625
$ws.Range("M89").Value = -50764.625
# Row 99
$ws.Range("H99").Value = 3030600.8
$ws.Range("I99").Value = 2377033
$ws.Range("J99").Value = 4010952.5
$ws.Range("K99").Value = 2377033
$ws.Range("L99").Value = 4010952.5
$ws.Range("M99").Value = -2375535
$ws.Range("N99").Value = -4013948.5
# Row 126
$ws.Range("H126").Value = 3030600.8
$ws.Range("I126").Value = 2377033
$ws.Range("J126").Value = 4010952.5
$ws.Range("K126").Value = 7131099
$ws.Range("L126").Value = 12032857.5
$ws.Range("M126").Value = -7128629
$ws.Range("N126").Value = -12037797.5
# Row 134
$ws.Range("H134").Value = 26321068
$ws.Range("I134").Value = 1732.4783
$ws.Range("J134").Value = 66677384
$ws.Range("K134").Value = 5197.4349
$ws.Range("L134").Value = 200032152
$ws.Range("M134").Value = -2662.4349
$ws.Range("N134").Value = -200037222

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 52.192307
$ws.Range("I2").Value = 48.916668
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 293.500008
$ws.Range("L2").Value = 330
$ws.Range("M2").Value = -180.500008
$ws.Range("N2").Value = -556
# Row 38
$ws.Range("H38").Value = 91.38461
$ws.Range("J38").Value = 215.75
$ws.Range("L38").Value = 647.25
$ws.Range("N38").Value = -1341.25
# Row 59
$ws.Range("H59").Value = 1839.6666
$ws.Range("I59").Value = 210
$ws.Range("J59").Value = 2654.5
$ws.Range("K59").Value = 630
$ws.Range("L59").Value = 7963.5
$ws.Range("M59").Value = -90
$ws.Range("N59").Value = -9043.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
# Row 40
$ws.Range("H40").Value = 18330
$ws.Range("I40").Value = 15000
$ws.Range("K40").Value = 15000
$ws.Range("M40").Value = -14849
# Row 107
$ws.Range("H107").Value = 1096.5
$ws.Range("I107").Value = 180
$ws.Range("K107").Value = 180
$ws.Range("M107").Value = 1740
# Row 122
$ws.Range("H122").Value = 1359554
$ws.Range("I122").Value = 1698892.5
$ws.Range("K122").Value = 5096677.5
$ws.Range("M122").Value = -5094227.5
# Cells removed entirely in the updated workbook
$ws.Range("N22").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 10000
$ws.Range("J13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("N13").Value = -10280
# Row 46
$ws.Range("H46").Value = 3221.625
$ws.Range("J46").Value = 3221.625
$ws.Range("L46").Value = 3221.625
$ws.Range("N46").Value = -3597.625
# Row 55
$ws.Range("H55").Value = 2249.25
$ws.Range("I55").Value = 1812.4286
$ws.Range("K55").Value = 1812.4286
$ws.Range("M55").Value = -1639.4286
# Row 122
$ws.Range("H122").Value = 20598258
$ws.Range("I122").Value = 35916480
$ws.Range("K122").Value = 107749440
$ws.Range("M122").Value = -107746990
# Row 132
$ws.Range("H132").Value = 2237125.5
$ws.Range("I132").Value = 2105.75
$ws.Range("J132").Value = 4025141.2
$ws.Range("K132").Value = 6317.25
$ws.Range("L132").Value = 12075423.6
$ws.Range("M132").Value = -3787.25
$ws.Range("N132").Value = -12080483.6
# Row 136
$ws.Range("H136").Value = 19391.479
$ws.Range("I136").Value = 16928.5
$ws.Range("K136").Value = 50785.5
$ws.Range("M136").Value = -48235.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 27509
$ws.Range("I62").Value = 22888.625
$ws.Range("J62").Value = 36749.75
$ws.Range("K62").Value = 22888.625
$ws.Range("L62").Value = 36749.75
$ws.Range("M62").Value = -22264.625
$ws.Range("N62").Value = -37997.75
# Row 65
$ws.Range("H65").Value = 27509
$ws.Range("I65").Value = 22888.625
$ws.Range("J65").Value = 36749.75
$ws.Range("K65").Value = 114443.125
$ws.Range("L65").Value = 183748.75
$ws.Range("M65").Value = -111323.125
$ws.Range("N65").Value = -189988.75
# Row 107
$ws.Range("H107").Value = 2285.4285
$ws.Range("I107").Value = 2499.6667
$ws.Range("K107").Value = 7499.000100000001
$ws.Range("M107").Value = -5579.000100000001
# Row 122
$ws.Range("H122").Value = 329834.75
$ws.Range("I122").Value = 421467.56
$ws.Range("K122").Value = 1264402.68
$ws.Range("M122").Value = -1261952.68
